# Add NOAA-20 platform -> VIIRS sensor mapping.
#
# The row is inserted at sheet row 163 (pushing the existing rows
# 163..202 down to 164..203), matching the hasSensor/hasPlatform
# mapping table layout already used throughout Sheet1
# (columns: A=Platform URI, B=relationship, C=Sensor URI, D=relationship).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Shift rows 163:202 down one row to make room for the new mapping.
$ws.Rows("163:163").Insert()

# Populate the newly inserted row with the NOAA-20 / VIIRS mapping.
$ws.Range("A163").Value = "plat_noaa_20"
$ws.Range("B163").Value = "hasSensor"
$ws.Range("C163").Value = "sens_viirs"
$ws.Range("D163").Value = "hasPlatform"
